# Raw and Clean Data from SSA for October 23-26
# Adds the new daily observations (2020-10-23 .. 2020-10-26) across the
# tracking sheets of the "bitacora historica" workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) out_vars : four new daily rows (146-149) for 10/23 .. 10/26
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("out_vars")

# Carry the formatting (date format on col A, General elsewhere) down from
# the last existing row before filling in the new values.
$ws.Range("A145:J145").Copy($ws.Range("A146:J149"))

$outVarsRows = @{
    146 = @(44127, 880775, 1058102, 330956, 88312, 22.990548096846528, 202495, 17422, 34076, 2269833)
    147 = @(44128, 886800, 1066646, 335143, 88743, 22.937302661253948, 203408, 17539, 34177, 2288589)
    148 = @(44129, 891160, 1072760, 331758, 88924, 22.897796130885588, 204056, 17575, 34244, 2295678)
    149 = @(44130, 895326, 1078072, 328231, 89171, 22.885183720789971, 204897, 17651, 34310, 2301629)
}
foreach ($r in 146..149) {
    $vals = $outVarsRows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $vals[$c]
    }
}
$ws.Range("A145").Select()

# ---------------------------------------------------------------------------
# 2) dates_dx : fill rows 142-144, add rows 145 and 147 (146 stays absent)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_dx")

$ws.Range("O141").Value2 = 5

$datesDxRows = @{
    142 = @(0,1,0,0,1,0,0,0,1,0,0,1,2,5)
    143 = @(0,1,0,0,1,0,0,0,1,0,0,1,2,5)
    144 = @(0,1,0,0,1,0,0,0,1,0,0,1,2,5)
}
foreach ($r in 142..144) {
    $vals = $datesDxRows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value2 = $vals[$c]
    }
}

# Row 145 (new) - 2020-10-23
$ws.Range("A144").Copy($ws.Range("A145"))
$ws.Range("A145").Value2 = 44127
$datesDxRow145 = @(0,1,0,0,1,0,0,0,1,0,0,1,2,5)
for ($c = 0; $c -lt $datesDxRow145.Length; $c++) {
    $ws.Cells.Item(145, $c + 2).Value2 = $datesDxRow145[$c]
}

# Row 147 (new, row 146 intentionally left blank) - 2020-10-25
$ws.Range("A144").Copy($ws.Range("A147"))
$ws.Range("A147").Value2 = 44129
$datesDxRow147 = @(0,1,0,0,1,0,0,0,1,0,0,1,2,5)
for ($c = 0; $c -lt $datesDxRow147.Length; $c++) {
    $ws.Cells.Item(147, $c + 2).Value2 = $datesDxRow147[$c]
}
$ws.Range("O147").Select()

# ---------------------------------------------------------------------------
# 3) dates_sx : add rows 142-145 and 147 (146 stays absent)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_sx")

$ws.Range("A141").Copy($ws.Range("A142:A143"))
$ws.Range("A142").Value2 = 44124
$ws.Range("A143").Value2 = 44125

$datesSxRows = @{
    142 = @(0,1,0,0,0,0,1,0,0,1,0,0,0)
    143 = @(0,1,0,0,0,0,1,0,0,1,0,0,0)
    144 = @(0,1,0,0,0,0,1,0,0,1,0,0,0)
}
foreach ($r in 142..144) {
    $vals = $datesSxRows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value2 = $vals[$c]
    }
}

$ws.Range("A144").Copy($ws.Range("A145"))
$ws.Range("A145").Value2 = 44127
$datesSxRow145 = @(0,1,0,0,0,0,1,0,0,1,0,0,0)
for ($c = 0; $c -lt $datesSxRow145.Length; $c++) {
    $ws.Cells.Item(145, $c + 2).Value2 = $datesSxRow145[$c]
}

$ws.Range("A144").Copy($ws.Range("A147"))
$ws.Range("A147").Value2 = 44129
$datesSxRow147 = @(0,1,0,0,0,0,1,0,0,1,0,0,0)
for ($c = 0; $c -lt $datesSxRow147.Length; $c++) {
    $ws.Cells.Item(147, $c + 2).Value2 = $datesSxRow147[$c]
}
$ws.Range("O147").Select()

# ---------------------------------------------------------------------------
# 4) dates_deaths : fill row 142, add rows 143-145 and 147 (146 absent)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_deaths")

$datesDeathsRow142 = @(0,0,0,1,1,1,0,2,1,2,1,2)
for ($c = 0; $c -lt $datesDeathsRow142.Length; $c++) {
    $ws.Cells.Item(142, $c + 2).Value2 = $datesDeathsRow142[$c]
}

$ws.Range("A142").Copy($ws.Range("A143"))
$ws.Range("A143").Value2 = 44125
$datesDeathsRow143 = @(0,0,0,1,1,1,0,2,1,2,1,2)
for ($c = 0; $c -lt $datesDeathsRow143.Length; $c++) {
    $ws.Cells.Item(143, $c + 2).Value2 = $datesDeathsRow143[$c]
}

$datesDeathsRow144 = @(0,0,0,1,1,1,0,2,1,2,1,2)
for ($c = 0; $c -lt $datesDeathsRow144.Length; $c++) {
    $ws.Cells.Item(144, $c + 2).Value2 = $datesDeathsRow144[$c]
}

$ws.Range("A144").Copy($ws.Range("A145"))
$ws.Range("A145").Value2 = 44127
$datesDeathsRow145 = @(0,0,0,1,1,1,0,2,1,2,1,2)
for ($c = 0; $c -lt $datesDeathsRow145.Length; $c++) {
    $ws.Cells.Item(145, $c + 2).Value2 = $datesDeathsRow145[$c]
}

$ws.Range("A144").Copy($ws.Range("A147"))
$ws.Range("A147").Value2 = 44129
$datesDeathsRow147 = @(0,0,0,1,1,1,0,2,1,2,1,2)
for ($c = 0; $c -lt $datesDeathsRow147.Length; $c++) {
    $ws.Cells.Item(147, $c + 2).Value2 = $datesDeathsRow147[$c]
}
$ws.Range("N147").Select()

# ---------------------------------------------------------------------------
# 5) control_obs : 4 new date columns EP:ES (10/23 .. 10/26)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("control_obs")

$ws.Range("EP1").Value2 = 44127
$ws.Range("EQ1").Value2 = 44128
$ws.Range("ER1").Value2 = 44129
$ws.Range("ES1").Value2 = 44130

$controlObsRows = @{
    2  = @(7609, $null, 7673, 7705)
    3  = @(7452, $null, 7516, 7548)
    4  = @(7452, $null, 7516, 7548)
    5  = @(7452, $null, 7516, 7548)
    6  = @(7452, $null, 7516, 7548)
    7  = @(6688, $null, 6752, 6784)
    8  = @(9483, $null, 9547, 9579)
    10 = @(285,  $null, 287,  288)
    11 = @(285,  $null, 287,  288)
    12 = @(285,  $null, 287,  288)
    13 = @(285,  $null, 287,  288)
    14 = @(285,  $null, 287,  288)
    15 = @(235,  $null, 236,  238)
    16 = @(297,  $null, 298,  300)
    18 = @(1752, $null, 1766, 1773)
}
$cols = @("EP", "EQ", "ER", "ES")
foreach ($r in $controlObsRows.Keys) {
    $vals = $controlObsRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Range($cols[$i] + $r).Value2 = $vals[$i]
        }
    }
}

$ws.Range("EO20").Copy($ws.Range("EP20:ES20"))
foreach ($col in $cols) {
    $ws.Range($col + "20").Formula = "=SUM(" + $col + "2:" + $col + "18)"
}

$ws.Range("EW14").Select()

Write-Output "done"
